$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update invoice header: HD number + date/time (this is a refund of the same
#        sale a bit later the same day, per the commit message) ---
$ws.Range("A6").Value = "Số HD: HD120424008"
$ws.Range("A7").Value = "Ngày giờ: 16:47:00 - 12/04/2024"

# --- 2. "Mắt biếc" line: quantity 1 -> 3, so its line total becomes 181,800 ---
$ws.Range("C13").Value = 3
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "'181,800"
$ws.Range("D13").NumberFormat = "General"

# --- 3. The "Kéo" line (old rows 16:17) is removed entirely (refund of that
#        item), shifting everything below it up by two rows ---
$ws.Rows("16:17").Delete()

# --- 4. Refresh the totals so they reflect only "Mắt biếc" (181,800) + "BatMan"
#        (87,360) = 269,160. After the delete above the summary block now sits at
#        rows 16-20: Tạm tính(16) / Giảm giá(17) / Tổng cộng(18) / Tiền khách đưa(19)
#        / Tiền trả lại khách(20). Giảm giá and Tiền trả lại khách stay "0". ---
foreach ($addr in @("D16","D18","D19")) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = "'269,160"
    $c.NumberFormat = "General"
}
